# Update cryptocurrency price list data (cryptos.xlsx)
# Refresh Price / Volume(1h) columns with latest scraped values.
# Two coin pairs changed rank order, so their Coin name / Link / Price / Volume
# cells are updated in place on the corresponding rows as well.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'28.197.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.60%  "

# Row 3
$ws.Range("D3").Value = "'1.805.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.86%  "

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.45%  "

# Row 5
$ws.Range("D5").Value = "'312.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.07%  "

# Row 6
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.42%  "

# Row 7
$ws.Range("D7").Value = "'0.5137"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("D8").Value = "'0.3967"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.26%  "

# Row 9
$ws.Range("D9").Value = "'0.07800"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.46%  "

# Row 10
$ws.Range("D10").Value = "'1.107"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.77%  "

# Row 11
$ws.Range("D11").Value = "'40.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.38%  "

# Row 12
$ws.Range("D12").Value = "'6.354"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.76%  "

# Row 13
$ws.Range("D13").Value = "'1.001"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.43%  "

# Row 14
$ws.Range("D14").Value = "'20.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.75%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'7.303"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.58%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'1.801.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.08%  "

# Row 17
$ws.Range("D17").Value = "'92.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.95%  "

# Row 18
$ws.Range("D18").Value = "'0.00001078"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.33%  "

# Row 19
$ws.Range("D19").Value = "'0.06569"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.89%  "

# Row 20
$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.36%  "

# Row 21
$ws.Range("D21").Value = "'17.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.32%  "

# Row 22
$ws.Range("D22").Value = "'5.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.44%  "

# Row 23
$ws.Range("D23").Value = "'28.227.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.65%  "

# Row 24
$ws.Range("D24").Value = "'11.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.74%  "

# Row 25
$ws.Range("D25").Value = "'2.213"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.34%  "

# Row 26
$ws.Range("D26").Value = "'160.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.67%  "

# Row 27
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'2.434"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.86%  "

# Row 28
$ws.Range("D28").Value = "'20.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.97%  "

# Row 29
$ws.Range("B29").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C29").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D29").Value = "'2.017.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.50%  "

# Row 30
$ws.Range("D30").Value = "'127.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.47%  "

# Row 31
$ws.Range("D31").Value = "'0.1094"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.06%  "

# Row 32
$ws.Range("D32").Value = "'1.057"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.17%  "

# Row 33
$ws.Range("D33").Value = "'3.654"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.78%  "

# Row 34
$ws.Range("D34").Value = "'5.563"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.15%  "

# Row 35
$ws.Range("D35").Value = "'0.07168"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.22%  "

# Row 36
$ws.Range("D36").Value = "'9.095"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.43%  "

# Row 37
$ws.Range("D37").Value = "'0.02348"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.42%  "

# Row 38
$ws.Range("D38").Value = "'0.2173"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.21%  "

# Row 39
$ws.Range("D39").Value = "'5.034"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.23%  "

# Row 40
$ws.Range("D40").Value = "'11.53"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.15%  "

# Row 41
$ws.Range("D41").Value = "'0.6159"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.50%  "

# Row 42
$ws.Range("D42").Value = "'1.001"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.42%  "

# Row 43
$ws.Range("D43").Value = "'1.153"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.19%  "

# Row 44
$ws.Range("D44").Value = "'13.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.79%  "

# Row 45
$ws.Range("D45").Value = "'1.309"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.66%  "

# Row 46
$ws.Range("D46").Value = "'0.5954"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.96%  "

# Row 47
$ws.Range("D47").Value = "'3.738"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.63%  "

# Row 48
$ws.Range("D48").Value = "'124.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.79%  "

# Row 49
$ws.Range("D49").Value = "'1.214"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.29%  "

# Row 50
$ws.Range("D50").Value = "'1.915"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.82%  "

# Row 51
$ws.Range("D51").Value = "'0.06796"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.50%  "
